$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text type for the Price (D) and Hora (G) columns, which are stored as
# text in the source data (e.g. "248.05", "8") rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "247.95"
$ws.Range("G2").Value = "9"

# Row 3
$ws.Range("D3").Value = "22.01"
$ws.Range("G3").Value = "9"

# Row 4
$ws.Range("D4").Value = "5.334"
$ws.Range("G4").Value = "9"

# Row 5
$ws.Range("D5").Value = "0.05634"
$ws.Range("G5").Value = "9"

# Row 6
$ws.Range("D6").Value = "3.427"
$ws.Range("G6").Value = "9"

# Row 7
$ws.Range("D7").Value = "6.375"
$ws.Range("G7").Value = "9"

# Row 8
$ws.Range("D8").Value = "0.8187"
$ws.Range("G8").Value = "9"

# Row 9
$ws.Range("D9").Value = "0.9374"
$ws.Range("G9").Value = "9"

# Row 10
$ws.Range("D10").Value = "0.1441"
$ws.Range("G10").Value = "9"

# Row 11
$ws.Range("D11").Value = "0.07514"
$ws.Range("G11").Value = "9"

# Row 12
$ws.Range("D12").Value = "0.03247"
$ws.Range("G12").Value = "9"

# Row 13
$ws.Range("D13").Value = "0.03085"
$ws.Range("G13").Value = "9"

# Row 14
$ws.Range("D14").Value = "0.09311"
$ws.Range("G14").Value = "9"

# Row 15
$ws.Range("D15").Value = "3.579"
$ws.Range("G15").Value = "9"

# Row 16
$ws.Range("G16").Value = "9"

# Row 17
$ws.Range("D17").Value = "0.04736"
$ws.Range("G17").Value = "9"

# Row 18
$ws.Range("D18").Value = "0.01149"
$ws.Range("E18").Value = "17OneONEBestin24h"
$ws.Range("G18").Value = "9"

# Row 19
$ws.Range("D19").Value = "0.006267"
$ws.Range("G19").Value = "9"

# Row 20
$ws.Range("D20").Value = "0.005062"
$ws.Range("G20").Value = "9"

# Row 21
$ws.Range("D21").Value = "0.001032"
$ws.Range("G21").Value = "9"

# Row 22
$ws.Range("G22").Value = "9"

# Row 23
$ws.Range("D23").Value = "3.760"
$ws.Range("G23").Value = "9"

# Row 24
$ws.Range("D24").Value = "2.147"
$ws.Range("G24").Value = "9"

# Row 25
$ws.Range("D25").Value = "0.3307"
$ws.Range("G25").Value = "9"

# Row 26
$ws.Range("G26").Value = "9"

# Row 27
$ws.Range("G27").Value = "9"

# Row 28
$ws.Range("G28").Value = "9"

# Row 29
$ws.Range("G29").Value = "9"

# Row 30
$ws.Range("G30").Value = "9"

# Row 31
$ws.Range("G31").Value = "9"

# Row 32
$ws.Range("G32").Value = "9"

# Row 33
$ws.Range("G33").Value = "9"

# Row 34
$ws.Range("G34").Value = "9"

# Row 35
$ws.Range("G35").Value = "9"

# Row 36
$ws.Range("G36").Value = "9"

# Row 37
$ws.Range("G37").Value = "9"

# Row 38
$ws.Range("G38").Value = "9"

# Row 39
$ws.Range("G39").Value = "9"

# Row 40
$ws.Range("D40").Value = "0.03961"
$ws.Range("G40").Value = "9"

# Row 41
$ws.Range("D41").Value = "0.006987"
$ws.Range("G41").Value = "9"

# Row 42
$ws.Range("D42").Value = "0.1065"
$ws.Range("G42").Value = "9"

# Row 43
$ws.Range("D43").Value = "0.003398"
$ws.Range("G43").Value = "9"

# Row 44
$ws.Range("D44").Value = "0.008826"
$ws.Range("G44").Value = "9"

# Row 45
$ws.Range("G45").Value = "9"

# Row 46
$ws.Range("G46").Value = "9"

# Row 47
$ws.Range("D47").Value = "0.0005498"
$ws.Range("G47").Value = "9"

# Row 48
$ws.Range("D48").Value = "0.7796"
$ws.Range("G48").Value = "9"

# Row 49
$ws.Range("D49").Value = "0.1769"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("G49").Value = "9"

# Row 50
$ws.Range("G50").Value = "9"

# Row 51
$ws.Range("G51").Value = "9"
